$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.652.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.895.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4920"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2943"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06695"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.889.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07339"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.137"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6652"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.642.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007860"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.149.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("B21").Value = "BinanceUSD"
$ws.Range("C21").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.294"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "189.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.210"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.540"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.924"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.477"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.389"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09133"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.014"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05232"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7377"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.100"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01818"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.692"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9175"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.059"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +28.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4410"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.921"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9964"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1380"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.522"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.070"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.50%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05854"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3954"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.46%  "
